$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-09-06 15:25:38"

$wsZhCn.Range("H2").Value = "2016-09-06 15:25:20"
$wsZhCn.Range("K2").Value = "2016-09-06 15:26:37"

$wsDeDe.Range("K2").Value = "2016-09-06 15:26:56"
